$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logindata")

$ws.Range("A2").Value = "locked_out_user"

$ws.Range("A3").Value = "problem_user"
$ws.Range("B3").Value = "secret_sauce"

$ws.Range("A4").Value = "performance_glitch_user"
$ws.Range("B4").Value = "secret_sauce"

$ws.Range("B9").Select()
